# "search needed multipliers from gain_setting file"
#
# - Rename the first two sheets (44AMU -> Inert, 46AMU -> Reactant).
# - Pull the "temperature" column (N) that already lives on Reactant (46AMU)
#   over onto 18AMU, since that sheet was still missing the gain-setting
#   temperature readings used to derive the needed multipliers.
# - Leave the selection on Reactant parked on the whole temperature column,
#   and finish with 18AMU active/selected at Q8 (where the lookup work
#   happens), which also drops the previously-active tab on 16AMU.

$wb = $excel.ActiveWorkbook

$wsInert    = $wb.Worksheets.Item(1)
$wsInert.Name = "Inert"

$wsReactant = $wb.Worksheets.Item(2)
$wsReactant.Name = "Reactant"

$ws18AMU = $wb.Worksheets.Item(3)

# Reactant has a 14th column ("temperature") that 18AMU is missing -
# copy it across, header included, row by row.
$lastRow = 182
$tempCol = 14

$wsReactant.Activate()
$wsReactant.Columns("N").Select()

for ($r = 1; $r -le $lastRow; $r++) {
    $ws18AMU.Cells.Item($r, $tempCol).Value = $wsReactant.Cells.Item($r, $tempCol).Value2
}

$ws18AMU.Activate()
$ws18AMU.Range("Q8").Select()
